# Terminal La Palmera de La Serena - Cebollín: add a new weekly data row.
#
# The underlying data rows for this market/category (rows 81-144, columns
# D "Fecha", J "Volumen", K "Precio minimo", L "Precio maximo",
# M "Precio promedio ponderado" and P "Precio $/Kg") get shifted down by
# one row - i.e. every row takes on the values that used to belong to the
# row above it - while a brand new observation is written into row 81 and
# the observation that used to sit in row 144 is pushed out into a new
# row 145 (copied verbatim, since none of its other columns change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 81
$lastRow  = 144
$newRow   = 145

# --- 1. Snapshot every column for the last existing row (144) so it can be
#        duplicated into the brand-new row 145 before anything is overwritten.
$lastRowVals = @{}
for ($col = 1; $col -le 18; $col++) {
    $lastRowVals[$col] = $ws.Cells.Item($lastRow, $col).Value2
}
$dateFormat = $ws.Cells.Item($lastRow, 4).NumberFormat

# --- 2. Snapshot the columns that shift (D, J, K, L, M, P) for every row in
#        the block, before any writes happen.
$dArr = @{}
$jArr = @{}
$kArr = @{}
$lArr = @{}
$mArr = @{}
$pArr = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dArr[$r] = $ws.Cells.Item($r, 4).Value2
    $jArr[$r] = $ws.Cells.Item($r, 10).Value2
    $kArr[$r] = $ws.Cells.Item($r, 11).Value2
    $lArr[$r] = $ws.Cells.Item($r, 12).Value2
    $mArr[$r] = $ws.Cells.Item($r, 13).Value2
    $pArr[$r] = $ws.Cells.Item($r, 16).Value2
}

# --- 3. Write the new row 145 as an exact copy of the old row 144.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $lastRowVals[$col]
}
$ws.Cells.Item($newRow, 4).NumberFormat = $dateFormat

# --- 4. Shift rows 144 down to 82: each row gets the shifting columns from
#        the row immediately above it (walk backwards so we never read a
#        cell that a previous loop iteration has already overwritten -
#        not strictly required since we snapshotted above, but keeps intent
#        clear).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $srcRow = $r - 1
    $ws.Cells.Item($r, 4).Value  = $dArr[$srcRow]
    $ws.Cells.Item($r, 10).Value = $jArr[$srcRow]
    $ws.Cells.Item($r, 11).Value = $kArr[$srcRow]
    $ws.Cells.Item($r, 12).Value = $lArr[$srcRow]
    $ws.Cells.Item($r, 13).Value = $mArr[$srcRow]
    $ws.Cells.Item($r, 16).Value = $pArr[$srcRow]
}

# --- 5. Row 81 becomes the new observation: new date + new volume, the
#        min/max/avg price and $/Kg columns are unchanged.
$ws.Cells.Item($firstRow, 4).Value  = 44554
$ws.Cells.Item($firstRow, 10).Value = 3000
